# Auto-generated edit script: apply numeric value updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 284.27274
$ws.Range("I4").Value = 284.27274
$ws.Range("K4").Value = 284.27274
$ws.Range("M4").Value = -170.27274
$ws.Range("H88").Value = 11572.429
$ws.Range("I88").Value = 6501.5
$ws.Range("K88").Value = 6501.5
$ws.Range("M88").Value = -6095.5
$ws.Range("H91").Value = 11572.429
$ws.Range("I91").Value = 6501.5
$ws.Range("K91").Value = 6501.5
$ws.Range("M91").Value = -5097.5
$ws.Range("H132").Value = 3045.651
$ws.Range("I132").Value = 1472.1818
$ws.Range("J132").Value = 13863.25
$ws.Range("K132").Value = 4416.5454
$ws.Range("L132").Value = 41589.75
$ws.Range("M132").Value = -1886.5454
$ws.Range("N132").Value = -46649.75
$ws.Range("H138").Value = 2129964.2
$ws.Range("I138").Value = 890.8108
$ws.Range("J138").Value = 3511994.5
$ws.Range("K138").Value = 2672.4324
$ws.Range("L138").Value = 10535983.5
$ws.Range("M138").Value = 2467.5676
$ws.Range("N138").Value = -10546263.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2525.9167
$ws.Range("I2").Value = 2613.875
$ws.Range("J2").Value = 2350
$ws.Range("K2").Value = 2613.875
$ws.Range("L2").Value = 2350
$ws.Range("M2").Value = -2500.875
$ws.Range("N2").Value = -2576
$ws.Range("H4").Value = 173
$ws.Range("I4").Value = 173
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 173
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -57
$ws.Range("H63").Value = 3207.0667
$ws.Range("I63").Value = 2901
$ws.Range("J63").Value = 3411.111
$ws.Range("K63").Value = 2901
$ws.Range("L63").Value = 3411.111
$ws.Range("M63").Value = -2215
$ws.Range("N63").Value = -4783.111
$ws.Range("H66").Value = 3207.0667
$ws.Range("I66").Value = 2901
$ws.Range("J66").Value = 3411.111
$ws.Range("K66").Value = 14505
$ws.Range("L66").Value = 17055.555
$ws.Range("M66").Value = -11073
$ws.Range("N66").Value = -23919.555
$ws.Range("H74").Value = 23208.064
$ws.Range("I74").Value = 30507.676
$ws.Range("J74").Value = 2525.8333
$ws.Range("K74").Value = 30507.676
$ws.Range("L74").Value = 2525.8333
$ws.Range("M74").Value = -29633.676
$ws.Range("N74").Value = -4273.8333
$ws.Range("H77").Value = 23208.064
$ws.Range("I77").Value = 30507.676
$ws.Range("J77").Value = 2525.8333
$ws.Range("K77").Value = 152538.38
$ws.Range("L77").Value = 12629.1665
$ws.Range("M77").Value = -148170.38
$ws.Range("N77").Value = -21365.1665
$ws.Range("H116").Value = 2525.9167
$ws.Range("I116").Value = 2613.875
$ws.Range("J116").Value = 2350
$ws.Range("K116").Value = 2613.875
$ws.Range("L116").Value = 2350
$ws.Range("M116").Value = -319.875
$ws.Range("N116").Value = -6938
$ws.Range("H135").Value = 31666.125
$ws.Range("J135").Value = 31666.125
$ws.Range("L135").Value = 31666.125
$ws.Range("N135").Value = -41806.125
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2525.9167
$ws.Range("I3").Value = 2613.875
$ws.Range("J3").Value = 2350
$ws.Range("K3").Value = 2613.875
$ws.Range("L3").Value = 2350
$ws.Range("M3").Value = -2499.875
$ws.Range("N3").Value = -2578
$ws.Range("H103").Value = 26249.75
$ws.Range("J103").Value = 26249.75
$ws.Range("L103").Value = 26249.75
$ws.Range("N103").Value = -28593.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 24326.572
$ws.Range("J28").Value = 24326.572
$ws.Range("L28").Value = 24326.572
$ws.Range("N28").Value = -24816.572
$ws.Range("H31").Value = 2092.1177
$ws.Range("I31").Value = 996.8889
$ws.Range("J31").Value = 2486.4
$ws.Range("K31").Value = 996.8889
$ws.Range("L31").Value = 2486.4
$ws.Range("M31").Value = -701.8889
$ws.Range("N31").Value = -3076.4
$ws.Range("H34").Value = 2092.1177
$ws.Range("I34").Value = 996.8889
$ws.Range("J34").Value = 2486.4
$ws.Range("K34").Value = 996.8889
$ws.Range("L34").Value = 2486.4
$ws.Range("M34").Value = -794.8889
$ws.Range("N34").Value = -2890.4
$ws.Range("H58").Value = 3828.2703
$ws.Range("I58").Value = 4958.115
$ws.Range("J58").Value = 1157.7273
$ws.Range("K58").Value = 4958.115
$ws.Range("L58").Value = 1157.7273
$ws.Range("M58").Value = -4755.115
$ws.Range("N58").Value = -1563.7273
$ws.Range("H105").Value = 913.4783
$ws.Range("I105").Value = 857.8946999999999
$ws.Range("J105").Value = 1177.5
$ws.Range("K105").Value = 857.8946999999999
$ws.Range("L105").Value = 1177.5
$ws.Range("M105").Value = 889.1053000000001
$ws.Range("N105").Value = -4671.5
$ws.Range("H136").Value = 3828.2703
$ws.Range("I136").Value = 4958.115
$ws.Range("J136").Value = 1157.7273
$ws.Range("K136").Value = 14874.345
$ws.Range("L136").Value = 3473.1819
$ws.Range("M136").Value = -12324.345
$ws.Range("N136").Value = -8573.1819
$ws.Range("H141").Value = 78831
$ws.Range("J141").Value = 78831
$ws.Range("L141").Value = 78831
$ws.Range("N141").Value = -89191

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 325886.75
$ws.Range("I2").Value = 784.1539
$ws.Range("J2").Value = 627767.7
$ws.Range("K2").Value = 4704.9234
$ws.Range("L2").Value = 3766606.2
$ws.Range("M2").Value = -4591.9234
$ws.Range("N2").Value = -3766832.2
$ws.Range("H131").Value = 921.1900000000001
$ws.Range("I131").Value = 916
$ws.Range("J131").Value = 921.24243
$ws.Range("K131").Value = 2748
$ws.Range("L131").Value = 2763.72729
$ws.Range("M131").Value = 2292
$ws.Range("N131").Value = -12843.72729

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 133.6
$ws.Range("I2").Value = 39.77778
$ws.Range("J2").Value = 274.33334
$ws.Range("K2").Value = 39.77778
$ws.Range("L2").Value = 274.33334
$ws.Range("M2").Value = 73.22221999999999
$ws.Range("N2").Value = -500.33334
$ws.Range("H57").Value = 16427.092
$ws.Range("J57").Value = 17919.8
$ws.Range("L57").Value = 17919.8
$ws.Range("N57").Value = -19559.8
$ws.Range("H101").Value = 31500
$ws.Range("J101").Value = 31500
$ws.Range("L101").Value = 31500
$ws.Range("N101").Value = -37990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 74890.16
$ws.Range("I122").Value = 17900.908
$ws.Range("J122").Value = 153250.38
$ws.Range("K122").Value = 53702.724
$ws.Range("L122").Value = 459751.14
$ws.Range("M122").Value = -51252.724
$ws.Range("N122").Value = -464651.14
$ws.Range("H126").Value = 1041.8125
$ws.Range("I126").Value = 1017.7778
$ws.Range("K126").Value = 3053.3334
$ws.Range("M126").Value = -583.3334

